# Swap the embedded drawing "name" identifiers for the two logo images
# that appear in the document's headers/footers:
#   - Pearson logo (footers):  image2.png -> image1.png
#   - BTec logo   (headers):   image1.jpg -> image2.jpg
#
# These names live on <wp:docPr name="..."/> and <pic:cNvPr name="..."/>
# inside header1.xml/header2.xml/footer1.xml/footer2.xml. The Word object
# model has no direct InlineShape.Name property, so we round-trip the
# package through Document.WordOpenXML and do a targeted text swap of the
# "name=" attribute values only (the "descr=" / alt-text values are left
# untouched).

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$pearsonOld = 'name="image2.png"'
$pearsonNew = 'name="image1.png"'
$btecOld    = 'name="image1.jpg"'
$btecNew    = 'name="image2.jpg"'

$xml = $xml.Replace($pearsonOld, $pearsonNew)
$xml = $xml.Replace($btecOld, $btecNew)

$d.WordOpenXML = $xml
